$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.742.63"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "3.304.39"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +1.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.46"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.303.48"
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.56"
$ws.Range("E11").Value = "  +3.72%  "
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.86"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").Value = "3.850.20"
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "3.306.81"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("D18").Value = "63.817.23"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.28"
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.738"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.99"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.94"
$ws.Range("E24").Value = "  +5.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.18"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +1.69%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.74"
$ws.Range("E32").Value = "  +4.42%  "
$ws.Range("E33").Value = "  -1.23%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.05"
$ws.Range("E36").Value = "  +2.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.47"
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("E38").Value = "  +5.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0399"
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("D40").Value = "3.116.13"
$ws.Range("E40").Value = "  +4.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "431.42"
$ws.Range("E41").Value = "  +2.02%  "
$ws.Range("E42").Value = "  +7.62%  "
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.264"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("E46").Value = "  +2.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.72"
$ws.Range("E47").Value = "  +8.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.35"
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "126.77"
$ws.Range("E50").Value = "  +4.32%  "
$ws.Range("E51").Value = "  +0.09%  "
